# Add the new "2e Projectleider" column (Y) to the "Overzicht" sheet, mirroring
# the existing header formatting used for the other header cells (e.g. column X).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overzicht")

# New header cell in column Y, row 1.
$ws.Range("Y1").Value = "2e Projectleider"

# Copy the formatting (fill/font/alignment) of the neighbouring header cell (X1)
# onto the new header cell so it matches the rest of the header row.
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)

# Match the column width used by the other "bestFit" header columns (~17 chars).
$ws.Columns.Item(25).ColumnWidth = 16.166666666666668

# Make the new header cell the active selection, as in the edited workbook.
$ws.Range("Y1").Select() | Out-Null
